# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
#
# All of the target cells hold plain TEXT (e.g. "1.005", "29.081.56",
# "  -0.52%  ") even though many of them look numeric. Assigning such a
# string straight to Range.Value makes Excel "helpfully" reinterpret it as
# a number (losing the original text form / trailing zeros), and forcing
# text via NumberFormat="@" or a leading apostrophe instead bakes a new
# quotePrefix/number-format style onto the cell that the source workbook
# never had. Routing the literal through a temporary ="..." text formula
# and then Copy + PasteSpecial(xlPasteValues=-4163) bakes down to a plain
# literal string value with no formula and no style change - exactly
# matching the original cell's shape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $r = $ws.Range($rangeAddr)
    $escaped = $text.Replace("""", """""")
    $r.Formula = "=""" + $escaped + """"
    $r.Copy()
    $r.PasteSpecial(-4163)
}

Set-TextValue 'D2' '29.081.56'
Set-TextValue 'E2' '  -0.52%  '
Set-TextValue 'D3' '1.836.16'
Set-TextValue 'E3' '  -0.47%  '
Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  +0.54%  '
Set-TextValue 'D5' '242.32'
Set-TextValue 'E5' '  -0.21%  '
Set-TextValue 'D6' '0.6248'
Set-TextValue 'E6' '  -5.77%  '
Set-TextValue 'E7' '  +0.28%  '
Set-TextValue 'D8' '0.07590'
Set-TextValue 'E8' '  +1.98%  '
Set-TextValue 'D10' '22.70'
Set-TextValue 'E10' '  -2.68%  '
Set-TextValue 'D11' '0.07750'
Set-TextValue 'E11' '  -0.19%  '
Set-TextValue 'D12' '1.842.27'
Set-TextValue 'E12' '  -0.16%  '
Set-TextValue 'D13' '4.958'
Set-TextValue 'E13' '  -1.25%  '
Set-TextValue 'D14' '0.6653'
Set-TextValue 'E14' '  -1.06%  '
Set-TextValue 'D15' '0.000009949'
Set-TextValue 'E15' '  +13.90%  '
Set-TextValue 'D16' '82.84'
Set-TextValue 'E16' '  -0.80%  '
Set-TextValue 'D17' '6.023'
Set-TextValue 'E17' '  -2.67%  '
Set-TextValue 'D18' '29.134.33'
Set-TextValue 'E18' '  -0.34%  '
Set-TextValue 'D19' '225.48'
Set-TextValue 'E19' '  -0.67%  '
Set-TextValue 'D20' '12.35'
Set-TextValue 'E20' '  -1.59%  '
Set-TextValue 'D21' '1.002'
Set-TextValue 'E21' '  +0.19%  '
Set-TextValue 'D22' '7.210'
Set-TextValue 'E22' '  +0.43%  '
Set-TextValue 'D23' '1.003'
Set-TextValue 'E23' '  +0.32%  '
Set-TextValue 'D24' '158.76'
Set-TextValue 'E24' '  +0.09%  '
Set-TextValue 'D25' '8.456'
Set-TextValue 'E25' '  -2.12%  '
Set-TextValue 'D26' '0.1367'
Set-TextValue 'E26' '  -2.87%  '
Set-TextValue 'D27' '17.91'
Set-TextValue 'E27' '  -0.75%  '
Set-TextValue 'D28' '1.495'
Set-TextValue 'E28' '  -0.86%  '
Set-TextValue 'D29' '4.075'
Set-TextValue 'E29' '  -1.56%  '
Set-TextValue 'D30' '4.032'
Set-TextValue 'E30' '  -0.61%  '
Set-TextValue 'D31' '1.199'
Set-TextValue 'E31' '  +0.72%  '
Set-TextValue 'D32' '0.05204'
Set-TextValue 'E32' '  -2.41%  '
Set-TextValue 'E33' '  -0.96%  '
Set-TextValue 'D34' '0.7366'
Set-TextValue 'E34' '  -1.43%  '
Set-TextValue 'D35' '1.144'
Set-TextValue 'E35' '  -1.16%  '
Set-TextValue 'D36' '2.706'
Set-TextValue 'E36' '  +1.85%  '
Set-TextValue 'D37' '1.255.59'
Set-TextValue 'E37' '  -4.50%  '
Set-TextValue 'D38' '2.769'
Set-TextValue 'E38' '  +0.39%  '
Set-TextValue 'D39' '0.01783'
Set-TextValue 'E39' '  -1.10%  '
Set-TextValue 'D40' '6.322'
Set-TextValue 'E40' '  -1.25%  '
Set-TextValue 'D41' '0.8978'
Set-TextValue 'E41' '  -0.70%  '
Set-TextValue 'D42' '1.003'
Set-TextValue 'E42' '  +0.35%  '
Set-TextValue 'D43' '101.41'
Set-TextValue 'E43' '  -1.91%  '
Set-TextValue 'B44' 'BabyDogeCoin'
Set-TextValue 'C44' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D44' '0.00000000125'
Set-TextValue 'E44' '  +2.36%  '
Set-TextValue 'B45' 'RocketPoolETH'
Set-TextValue 'C45' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D45' '1.978.04'
Set-TextValue 'E45' '  -0.91%  '
Set-TextValue 'D46' '64.20'
Set-TextValue 'E46' '  -1.91%  '
Set-TextValue 'D47' '0.5122'
Set-TextValue 'E47' '  -0.45%  '
Set-TextValue 'D48' '0.4015'
Set-TextValue 'E48' '  -0.23%  '
Set-TextValue 'D49' '8.848'
Set-TextValue 'E49' '  +0.94%  '
Set-TextValue 'D50' '0.05757'
Set-TextValue 'D51' '1.640'
Set-TextValue 'E51' '  -6.48%  '
